$d = $word.ActiveDocument

# Remove the obsolete "<[Questions]Exclude.IMG>" placeholder run from the
# Questions/body cell (now that a real stamp image + Templater test cover
# this case, the Exclude.IMG merge-field placeholder text is no longer
# needed in the template).
$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute("<[Questions]Exclude.IMG>", $true, $true, $false, $false, $false, $true, 1, $false, "", 2)
